# Updates the "cryptos" price/volume table (columns D = Price, E = Volume(1h))
# to the latest scraped snapshot. Values are stored as plain text in the
# workbook (mixed formats like "27.545.21" alongside plain decimals), so for
# any new Price value that Excel would otherwise auto-parse as a number we
# briefly force text entry via NumberFormat "@" and then reset the cell back
# to the "Normal" style afterwards so no stray formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.545.21'
$ws.Range("E2").Value = '  -1.41%  '

$ws.Range("D3").Value = '1.845.63'
$ws.Range("E3").Value = '  -2.18%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.006'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -1.11%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '333.43'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.71%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.003'
$ws.Range("D6").Style = "Normal"

$ws.Range("E7").Value = '  -0.99%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3856'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.38%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '46.30'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.39%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07919'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.54%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9949'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.83%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '21.51'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.10%  '

$ws.Range("E13").Value = '  -1.80%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.934'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.41%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.125'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.04%  '

$ws.Range("E16").Value = '  -1.27%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '89.09'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.81%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06669'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.73%  '

$ws.Range("E19").Value = '  -0.97%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.08'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.26%  '

$ws.Range("E21").Value = '  -1.14%  '

$ws.Range("D22").Value = '27.553.22'
$ws.Range("E22").Value = '  -1.38%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.387'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.64%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.90'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.32%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '158.22'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.92%  '

$ws.Range("E27").Value = '  -2.42%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.108'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.92%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.414'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.25%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '120.06'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.92%  '

$ws.Range("E31").Value = '  +2.17%  '

$ws.Range("E32").Value = '  -1.53%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.587'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.81%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.290'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.90%  '

$ws.Range("E35").Value = '  -0.62%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06035'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.35%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02229'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.84%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.305'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.81%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.182'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.16%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5891'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.69%  '

$ws.Range("E41").Value = '  -1.43%  '

$ws.Range("E42").Value = '  +0.55%  '

$ws.Range("E43").Value = '  -2.18%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5587'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.43%  '

$ws.Range("E45").Value = '  -0.29%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.908'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.26%  '

$ws.Range("E47").Value = '  -2.36%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '110.99'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.64%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.054'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.23%  '

$ws.Range("E50").Value = '  -1.51%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '70.19'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.84%  '
